$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the two previously-empty rows (5 and 6) with the new i18n entries.
$ws.Range("A5").Value = "homePage.meta.title"
$ws.Range("B5").Value = "首页"
$ws.Range("C5").Value = "Home"

$ws.Range("A6").Value = "homePage.meta.description"
$ws.Range("B6").Value = "云极客工具，励志做最轻盈最好用的在线工具。以工匠精神打造功能丰富的在线工具，无需下载即可免费使用"

# Re-sort the data block (rows 1-6, header included) alphabetically by
# column A, same as the author did via Data > Sort in Excel (the header
# row stays in place because Header is set to "yes").
$sortRange = $ws.Range("A1:E6")
$keyRange = $ws.Range("A1:A6")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange, 0, 1, 0, 0)
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.Apply()

# Restore the final selection left by the author.
$ws.Range("E15").Select()

$wb.Save()
